$d = $word.ActiveDocument

# 1) "fabricante " -> "fornecedor " in the body text
#    ("Homologacao de fabricante de produto nacionalizado" -> "... fornecedor ...")
$d.Content.Find.Execute("fabricante ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "fornecedor ", 2)

# 2) Update the cached PAGE field result shown in the header from "3" to "2"
$header = $d.Sections(1).Headers(1)
$header.Range.Find.Execute("3", $true, $false, $false, $false, $false,
                            $true, 1, $false, "2", 2)
